$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column AO ("24-jul") with its values
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last header cell (AN1) onto the new header
# cell (AO1) so it keeps the bold / centered / bordered header style.
$ws1.Range("AN1").Copy()
$ws1.Range("AO1").PasteSpecial(-4122)
$ws1.Range("AO1").Value = "24-jul"

$ws1.Range("AO2").Value = 96.23999999999999
$ws1.Range("AO3").Value = 87.36
$ws1.Range("AO4").Value = 81.88
$ws1.Range("AO5").Value = 71.79000000000001
$ws1.Range("AO6").Value = 63.27
$ws1.Range("AO7").Value = 80.90000000000001
$ws1.Range("AO8").Value = 85.62
$ws1.Range("AO9").Value = 90.44
$ws1.Range("AO10").Value = 103.04
$ws1.Range("AO11").Value = 90.83
$ws1.Range("AO12").Value = 70.90000000000001
$ws1.Range("AO13").Value = 65.64
$ws1.Range("AO14").Value = 71.05
$ws1.Range("AO15").Value = 34.04
$ws1.Range("AO16").Value = 28
$ws1.Range("AO17").Value = 37.02
$ws1.Range("AO18").Value = 26.53
$ws1.Range("AO19").Value = 41.42
$ws1.Range("AO20").Value = 74.72
$ws1.Range("AO21").Value = 92.29000000000001
$ws1.Range("AO22").Value = 87.93000000000001
$ws1.Range("AO23").Value = 90.92
$ws1.Range("AO24").Value = 109.19
$ws1.Range("AO25").Value = 102.06

# ---------------------------------------------------------------------
# Sheet "Gaz": append a new row 38 for 2025-07-22
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")
# Force the date text to be stored as plain text (matching the other
# date cells in column A) instead of being auto-converted to a date
# serial number, then drop the temporary text format so no stray
# number formatting is left behind on the cell.
$ws2.Range("A38").NumberFormat = "@"
$ws2.Range("A38").Value = "2025-07-22"
$ws2.Range("A38").Style = "Normal"
$ws2.Range("B38").Value = 32.625

# ---------------------------------------------------------------------
# Sheet "CO2": append a new row 38 for 2025-07-22
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A38").NumberFormat = "@"
$ws3.Range("A38").Value = "2025-07-22"
$ws3.Range("A38").Style = "Normal"
$ws3.Range("B38").Value = 68.25
